$wb = $excel.ActiveWorkbook

# --- Basket_Course_Allocations: simplify allocated rooms (C9:C17) ---
$wsBCA = $wb.Worksheets.Item("Basket_Course_Allocations")
$wsBCA.Range("C9").Value = "C001"
$wsBCA.Range("C10").Value = "C101"
$wsBCA.Range("C11").Value = "C104"
$wsBCA.Range("C12").Value = "C202"
$wsBCA.Range("C13").Value = "C204"
$wsBCA.Range("C14").Value = "C102"
$wsBCA.Range("C15").Value = "C104"
$wsBCA.Range("C16").Value = "C202"
$wsBCA.Range("C17").Value = "C203"

# --- Executive_Summary: update generation timestamp (C3) ---
$wsES = $wb.Worksheets.Item("Executive_Summary")
$wsES.Range("C3").Value = "2026-01-26 12:46"

# --- Classroom_Allocation: re-balance room assignments (stricter pre/post-mid separation) ---
$wsCA = $wb.Worksheets.Item("Classroom_Allocation")
$wsCA.Range("G5").Value = "large classroom"
$wsCA.Range("H5").NumberFormat = "@"
$wsCA.Range("H5").Value = "120"
$wsCA.Range("I5").ClearContents()
$wsCA.Range("M5").Value = "C001"
$wsCA.Range("M6").Value = "C101"
$wsCA.Range("I7").Value = "Projector"
$wsCA.Range("M7").Value = "C104"
$wsCA.Range("I8").Value = "Projector"
$wsCA.Range("M8").Value = "C202"
$wsCA.Range("M9").Value = "C204"
$wsCA.Range("G17").Value = "classroom"
$wsCA.Range("H17").NumberFormat = "@"
$wsCA.Range("H17").Value = "96"
$wsCA.Range("I17").Value = "Projector"
$wsCA.Range("M17").Value = "C102"
$wsCA.Range("G18").Value = "classroom"
$wsCA.Range("H18").NumberFormat = "@"
$wsCA.Range("H18").Value = "96"
$wsCA.Range("I18").Value = "Projector"
$wsCA.Range("M18").Value = "C104"
$wsCA.Range("G21").Value = "large classroom"
$wsCA.Range("H21").NumberFormat = "@"
$wsCA.Range("H21").Value = "120"
$wsCA.Range("I21").ClearContents()
$wsCA.Range("M21").Value = "C001"
$wsCA.Range("M22").Value = "C101"
$wsCA.Range("I23").Value = "Projector"
$wsCA.Range("M23").Value = "C104"
$wsCA.Range("I24").Value = "Projector"
$wsCA.Range("M24").Value = "C202"
$wsCA.Range("M25").Value = "C204"
$wsCA.Range("G37").Value = "large classroom"
$wsCA.Range("H37").NumberFormat = "@"
$wsCA.Range("H37").Value = "120"
$wsCA.Range("I37").ClearContents()
$wsCA.Range("M37").Value = "C001"
$wsCA.Range("G38").Value = "classroom"
$wsCA.Range("H38").NumberFormat = "@"
$wsCA.Range("H38").Value = "96"
$wsCA.Range("I38").Value = "Projector"
$wsCA.Range("M38").Value = "C101"
$wsCA.Range("G39").Value = "classroom"
$wsCA.Range("H39").NumberFormat = "@"
$wsCA.Range("H39").Value = "96"
$wsCA.Range("M39").Value = "C104"
$wsCA.Range("M40").Value = "C202"
$wsCA.Range("M41").Value = "C204"
$wsCA.Range("G42").Value = "classroom"
$wsCA.Range("H42").NumberFormat = "@"
$wsCA.Range("H42").Value = "96"
$wsCA.Range("I42").Value = "Projector"
$wsCA.Range("M42").Value = "C102"
$wsCA.Range("G43").Value = "classroom"
$wsCA.Range("H43").NumberFormat = "@"
$wsCA.Range("H43").Value = "96"
$wsCA.Range("I43").Value = "Projector"
$wsCA.Range("M43").Value = "C104"
$wsCA.Range("G46").Value = "classroom"
$wsCA.Range("H46").NumberFormat = "@"
$wsCA.Range("H46").Value = "96"
$wsCA.Range("I46").Value = "Projector"
$wsCA.Range("M46").Value = "C102"
$wsCA.Range("G47").Value = "classroom"
$wsCA.Range("H47").NumberFormat = "@"
$wsCA.Range("H47").Value = "96"
$wsCA.Range("I47").Value = "Projector"
$wsCA.Range("M47").Value = "C104"
$wsCA.Range("G48").Value = "classroom"
$wsCA.Range("H48").NumberFormat = "@"
$wsCA.Range("H48").Value = "96"
$wsCA.Range("M48").Value = "C202"
$wsCA.Range("I49").Value = "TV"
$wsCA.Range("M49").Value = "C203"
